$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.763.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.508.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.50%  "

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.501.56"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.604"
$ws.Range("D8").Style = "Normal"

$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.656"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.142"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -10.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -12.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.071.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.59%  "

$ws.Range("E16").Value = "  -1.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.513.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "65.614.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.85%  "

$ws.Range("E21").Value = "  -8.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "386.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "617.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "62.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.57%  "

$ws.Range("E36").Value = "  -9.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "40.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.40%  "

$ws.Range("E38").Value = "  +0.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.387"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0733"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.129"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.950.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0397"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.40%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.10%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.128"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.39%  "

